# Added Mastery Slides 13, 14 and 15
# - H18, H19, H20 (WSQ 13/14/15 rows) get a "Slides Folder" marker.
# - H29, H30, H31 and H34 get a "Video" marker (H34 previously said
#   "Done in project").
# - The sheet view scrolls/selects near the newly-edited rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marks")
$ws.Activate()

$ws.Range("H18").Value = "Slides Folder"
$ws.Range("H19").Value = "Slides Folder"
$ws.Range("H20").Value = "Slides Folder"

$ws.Range("H29").Value = "Video"
$ws.Range("H30").Value = "Video"
$ws.Range("H31").Value = "Video"
$ws.Range("H34").Value = "Video"

# Update the view to match where the edits were made.
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E18").Select()
